$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 7.392889022827148
$ws.Range("B1").Value = 5.523151397705078
$ws.Range("C1").Value = 4.558609485626221
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 3.286228895187378
